# Update column F (dSF) values per repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = -2
    7  = -5
    8  = -1
    9  = -3
    10 = 0
    11 = -14
    13 = -5
    14 = -7
    15 = 1
    16 = 9
    18 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
